$d = $word.ActiveDocument
$RSQUO = [char]0x2019

# ---------------------------------------------------------------------------
# Paragraph 2: "Here is a picture ..." -> "Here's a picture ... will be. ..."
# Split into two Find/Replace calls so we don't touch the hyperlink run in
# the middle of the paragraph. The second call starts its search just after
# the closing ")" so the existing (non-hyperlink) run formatting is kept
# instead of inheriting the hyperlink's character style.
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2).Range
$p2.Find.Execute(
    "Here is a picture of how we think the system should be. We have a web-site (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ("Here" + $RSQUO + "s a picture of how we think the system will be. We have a web-site ("),
    2)

$p2b = $d.Paragraphs.Item(2).Range
$p2b.Find.Execute(
    " . Here users can download the program, and register a user so they can log in. The username and password been saved in a SQL database. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". Here users can download the program and register as a user so they can log in. The username and password is saved in a SQL database. ",
    2)

# ---------------------------------------------------------------------------
# Paragraph 3: fix grammar / verb agreement, drop the "wrong" sentence (it
# moves out of this paragraph entirely), and drop the trailing _GoBack
# bookmark (it will be re-added later in the new final paragraph).
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3).Range
$p3.Find.Execute(
    "The user start the freele.jar file, and get a message about to log in. The program been connected to the server. The server is connected to the SQL database and check if the password and username is right. If its right you will be sent to the chat rom.  If its wrong you will be sent to the website and can restore it there.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ("The user starts the freele.jar file and gets a message about to log in. The program is connected to the server. The server is connected to the SQL database and checks if the password and username is correct. If it" + $RSQUO + "s right you will be sent to the chat rom. "),
    2)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Paragraph 4: reword the chat-room description.
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4).Range
$p4.Find.Execute(
    "That is a common room for everyone who is signed in.  if you want to chat private with someone you can click on one the username and it will pop up a new window where you can chat private. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In the first chat room the user goes to is a public one where everyone who is signed inn can chat.  If you want to chat private with someone you can click on one of the username and it will pop up a new window where you can chat privately. ",
    2)

# ---------------------------------------------------------------------------
# New final paragraph describing the project stage, containing the
# relocated _GoBack bookmark.
# ---------------------------------------------------------------------------
$p4after = $d.Paragraphs.Item(4).Range
$p4after.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = "This is the first stage of the project, after everything is working with the chat we will add encryption and decryption and a few other function for the user."

$bmFind = $p5.Range.Duplicate
$bmFind.Find.Execute("and a ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmPoint = $d.Range($bmFind.End, $bmFind.End)
$d.Bookmarks.Add("_GoBack", $bmPoint)
